$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.352.99'
$ws.Range("E2").Value = '  -0.21%  '
$ws.Range("D3").Value = '1.850.41'
$ws.Range("E3").Value = '  +0.15%  '
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '0.9965'
$cell.Style = "Normal"
$ws.Range("E4").Value = '  -0.33%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '240.72'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +0.01%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.6350'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +1.12%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.9982'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -0.23%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.07572'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -1.41%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.2923'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +0.15%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '24.50'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -0.99%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.07737'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +0.01%  '
$ws.Range("D12").Value = '1.850.30'
$ws.Range("E12").Value = '  +0.65%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '5.022'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -0.08%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '0.6834'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("E15").Value = '  -2.96%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '83.26'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -0.24%  '
$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '6.139'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -0.50%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '29.382.80'
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '230.11'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +1.01%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '12.38'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -0.26%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '0.9979'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '7.471'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +0.89%  '
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '0.9984'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("B24").Value = 'Monero'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '158.87'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +0.78%  '
$ws.Range("B25").Value = 'Stellar'
$ws.Range("C25").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '0.1400'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +1.79%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '8.456'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +0.82%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '17.66'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '1.416'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +5.39%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '1.476'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +0.80%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '0.05689'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '4.142'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +0.69%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '4.055'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +0.73%  '
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '1.828'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -0.66%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.156'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -0.42%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.7004'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -1.14%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '2.579'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '1.248.16'
$ws.Range("E37").Value = '  +1.74%  '
$ws.Range("E38").Value = '  +2.03%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '2.727'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -1.79%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '6.523'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -0.20%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.9039'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.9977'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -0.27%  '
$ws.Range("B43").Value = 'RocketPoolETH'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D43").Value = '2.013.01'
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '102.03'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +0.32%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '65.96'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '7.145'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.1170'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +2.14%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '9.058'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("B49").Value = 'TheSandbox'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '0.3965'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -1.23%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.00000000115'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -4.52%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '1.676'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +0.25%  '
